$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet, add the new one right after it -----------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "first"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "second"

# --- Populate "second" with the same kind of age/gender/income rows ------
$data = @(
    @(24, "M", 2000),
    @(35, "F", 3100),
    @(28, "F", 3800),
    @(21, "F", 2800),
    @(31, "M", 3500)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
}

# --- Restore / set the on-screen selections for each sheet ---------------
[void]$ws1.Range("F7").Select()
[void]$ws2.Range("G12").Select()
